$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Remove the hyperlinks first -- the engine doesn't re-target
# hyperlink ranges when rows shift, so we rebuild them from scratch
# at the very end once everything else is in its final place.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# Drop the four events that are being retired from the schedule
# (delete bottom-up so the remaining row numbers don't shift under us):
#   row 9 - "Introducation to UNIX - Part II "
#   row 7 - "Customizing Your Graphs using GraphPad Prism 8 - Part II"
#   row 6 - "From a GenBank ID to a Phylogenetic Tree in MEGAX"
#   row 4 - "Preparing and Submitting Protein Structures to the NIAID 3D Printing Service"
# ------------------------------------------------------------------
$ws.Rows(9).Delete()
$ws.Rows(7).Delete()
$ws.Rows(6).Delete()
$ws.Rows(4).Delete()

# After the deletions, the two surviving NIAID events ("Building Shiny
# Apps" and "Molecular Visualization with Chimera - Part II") have
# shifted up to rows 4 and 5. Row 6 is now free for the new event.

# ------------------------------------------------------------------
# Add the new "Studying the Microbiome Using the Nephele Web Platform"
# webinar as row 6.
# ------------------------------------------------------------------
$ws.Range("A6").Value = "NIAID BioIT listserv"
$ws.Range("B6").Value = "Studying the Microbiome Using the Nephele Web Platform"
$ws.Range("C6").Value = "4/21/2020"
$ws.Range("D6").Value = "4/21/2020"
$ws.Range("E6").Value = 1587474000
$ws.Range("F6").Value = "GoToWebinar (1:00 – 2:00 pm))"
$ws.Range("G6").Value = "https://attendee.gotowebinar.com/register/4895163611488872973"
$ws.Range("H6").Value = "Webinar"
$ws.Range("I6").Value = "yes"
$ws.Range("J6").Value = "This webinar will demonstrate how to process and analyze a 16S microbiome dataset as well as a shotgun metagenomics dataset using the pipelines available in the Nephele platform."
$ws.Range("K6").Value = "microbiome,analysis,cloud computing"

# ------------------------------------------------------------------
# Re-create the two hyperlinks at their final locations:
#   G4 -> Building Shiny Apps registration link
#   G6 -> Studying the Microbiome Using the Nephele Web Platform registration link
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G4"), "https://attendee.gotowebinar.com/register/7582329347919390989")
$ws.Hyperlinks.Add($ws.Range("G6"), "https://attendee.gotowebinar.com/register/4895163611488872973")

# ------------------------------------------------------------------
# Match the saved selection state.
# ------------------------------------------------------------------
$ws.Range("K7").Select()
